$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: drop the "turnover " prefix from these three headers
$ws.Range("F1").Value = "gross amount"
$ws.Range("G1").Value = "net amount"
$ws.Range("I1").Value = "purchase count"
